$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (date, TCM, TCM-5, TCM-X)
$rows = @(
    @("10-09-2021", 115.58, 189.08, 102.08),
    @("13-09-2021", 115.48, 188.6,  102),
    @("14-09-2021", 114.61, 187.17, 101.24),
    @("15-09-2021", 114.56, 187.11, 101.21),
    @("16-09-2021", 114.38, 186.76, 101.05),
    @("20-09-2021", 113.92, 186.08, 100.59)
)

$startRow = 176
$r = $startRow
foreach ($row in $rows) {
    # Column A holds a date-like string; force Text format first so Excel
    # stores it as a plain shared string instead of auto-converting it to
    # a date serial number, then revert the cell style back to Normal so
    # no extra formatting is left behind on the cell.
    $cellA = $ws.Range("A" + $r)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]

    $r++
}
